# Update the "grouped matches" sheet (BDRC ID / 84000 ID grouping) with
# freshly compiled attributions. The (B,C) pairs are reassigned across
# rows 2-60 (values get shuffled to new rows as the script recompiled
# the grouping), and the row-60 entry for P00KG07267 drops the
# 'eft:sarvajnadeva' alias, keeping only 'eft:sarvanyadeva'.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("grouped matches")

$data = @(
    @(2, 'P8266', '{''eft:dharmatasila'', ''eft:ch-nyi-tsultrim''}'),
    @(3, 'P8217', '{''eft:t-jnanagarbha'', ''eft:jnanagarbha''}'),
    @(4, 'P8273', '{''eft:rin-chen-tsho'', ''eft:rinchen-tso''}'),
    @(5, 'P8213', '{''eft:vidyakarasimha''}'),
    @(6, 'P2548', '{''eft:prajnavarman'', ''eft:prajnavarma''}'),
    @(7, 'P8265', '{''eft:ratnaraksita''}'),
    @(8, 'P2956', '{''eft:krsnapandita''}'),
    @(9, 'P0TMP098', '{''eft:jinavara''}'),
    @(10, 'P2637', '{''eft:trakpa-gyaltsen''}'),
    @(11, 'P8245', '{''eft:buddhakaravarma''}'),
    @(12, 'P0TMP080', '{''eft:hwa-shang-zab-mo''}'),
    @(13, 'https://lod.dila.edu.tw/resource.php?id=A000089', '{''eft:siladharma''}'),
    @(14, 'P8261', '{''eft:munivarma'', ''eft:munivarman''}'),
    @(15, 'P3379', '{''eft:dipamkarasrijnana'', ''eft:dipamkara-srijnana''}'),
    @(16, 'P8268', '{''eft:buddhaprabha''}'),
    @(17, 'P8182', '{''eft:ska-ba-dpal-brtsegs'', ''eft:dpal-brtsegs'', ''eft:paltsek'', ''eft:kawa-paltsek-under-the-name-paltsek-raksita-'', ''eft:ban-de-dpal-brtsegs''}'),
    @(18, 'P4255', '{''eft:t-jnanagarbha'', ''eft:yesh-nyingpo'', ''eft:ye-shes-snying-po''}'),
    @(19, 'P8222', '{''eft:jnanasiddhi''}'),
    @(20, 'P8183', '{''eft:klu-i-rgyal-mtshan'', ''eft:cog-ro-klu-i-rgyal-mtshan''}'),
    @(21, 'P8267', '{''eft:vijayasila''}'),
    @(22, '?', '{''eft:sakyasena''}'),
    @(23, 'P3214 ', '{''eft:danasila''}'),
    @(24, 'P8210', '{''eft:danasila''}'),
    @(25, 'P8228', '{''eft:surendrabodhi''}'),
    @(26, 'P0TMPT007', '{''eft:rnam-par-mi-rtog-pa''}'),
    @(27, 'P8222 ', '{''eft:jnanasidhi''}'),
    @(28, 'P5651', '{''eft:pa-tshab-nyi-ma-grags''}'),
    @(29, 'P8269', '{''eft:dgon-gling-rma''}'),
    @(30, 'P1KG8854', '{''eft:silendrabodhi'', ''eft:surendrabodhi''}'),
    @(31, 'P8219', '{''eft:visuddhasimha''}'),
    @(32, 'P8220', '{''eft:devacandra''}'),
    @(33, 'P1KG8854 ', '{''eft:srilendrabodhi''}'),
    @(34, 'P4CZ16819', '{''eft:sakyaprabha''}'),
    @(35, 'P0TMP104', '{''eft:punyasambhava''}'),
    @(36, 'P8260', '{''eft:dpal-dbyangs''}'),
    @(37, 'P3285 ', '{''eft:sakya-yesh-''}'),
    @(38, 'P753', '{''eft:rin-chen-bzang-po''}'),
    @(39, 'P4CZ15137', '{''eft:kumarakalasa''}'),
    @(40, 'P3709 ', '{''eft:phakpa-sherab''}'),
    @(41, 'P4258', '{''eft:dpal-byor''}'),
    @(42, 'P8211', '{''eft:vidyakaraprabha''}'),
    @(43, 'P8249', '{''eft:dharmakara''}'),
    @(44, 'P8205 ', '{''eft:band-yesh-d-''}'),
    @(45, 'P8213 ', '{''eft:t-vidyakarasimha''}'),
    @(46, 'P4242', '{''eft:sherab-lekpa''}'),
    @(47, 'P4CZ16780 ', '{''eft:manjusrigarbha''}'),
    @(48, 'P8151', '{''eft:gayadhara''}'),
    @(49, 'P8206', '{''eft:celu''}'),
    @(50, 'P4263', '{''eft:dge-ba-dpal''}'),
    @(51, 'P8171', '{''eft:dharmasribhadra''}'),
    @(52, 'P3456', '{''eft:tshul-khrims-rgyal-ba''}'),
    @(53, 'P8093', '{''eft:kamalagupta''}'),
    @(54, 'P4259', '{''eft:dpal-gyi-lhun-po'', ''eft:palgyi-lh-npo'', ''eft:ban-de-dpal-gyi-lhun-po''}'),
    @(55, 'P0RK8', '{''eft:dharmapala''}'),
    @(56, 'P8263', '{''eft:leki-d-''}'),
    @(57, 'P0TMP092', '{''eft:anandasri-s-''}'),
    @(58, 'P8209', '{''eft:jinamitra'', ''eft:dzi-na-mi-tra-k-'', ''eft:jinamitra-k-''}'),
    @(59, 'P8205', '{''eft:yesh-d-'', ''eft:zhang-yesh-d-'', ''eft:yesh-d-ye-shes-sde-'', ''eft:band-yesh-de'', ''eft:ye-shes-sde''}'),
    @(60, 'P00KG07267', '{''eft:sarvanyadeva''}')
)

foreach ($entry in $data) {
    $row = $entry[0]
    $bdrcId = $entry[1]
    $ftId = $entry[2]
    $ws.Cells.Item($row, 2).Value = $bdrcId
    $ws.Cells.Item($row, 3).Value = $ftId
}
